# [Anmol Singh] Interface implements in all classes
#
# The underlying app (re)hashed every stored user Password value (the
# hashes are Base64-encoded digests of the "Original Password" column)
# and, as a side effect of saving from the live application, the DOB
# column's stored serial-date/time got a fresh "now" fractional value.
#
# This script reproduces those cell-level changes on the Customer,
# Driver and Owner sheets.

$wb = $excel.ActiveWorkbook

# New DOB serial value shared by every touched row.
$newDob = 28430.958761064816

# ---------------------------------------------------------------
# Customer sheet - Password column I, rows 2-8
# ---------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("Customer")

$wsCustomer.Range("I2").Value = "29k5k0C0rcGH3fcDqVCFEMFhsWNrCntUheu2eLSs/hU="
$wsCustomer.Range("I3").Value = "N/aKBnVkBI6GPiqP6w/1xuXlCw78NgTOWUbZI4WqVkM="
$wsCustomer.Range("I4").Value = "uWj5Jk72Xh0DprjFlAMhbA77L5Dunisg3aFCdGf5pww="
$wsCustomer.Range("I5").Value = "S0zNufK4i2mf+hQ9yadrX1W3bE0Qu8sZpSFiubOLzuc="
$wsCustomer.Range("I6").Value = "DpHYvZr2vgH0z6LEStjANi6OEHd1SfAbiFD6d/MPnRU="
$wsCustomer.Range("I7").Value = "X2U1uaF/gu0Hfq/m92/wY31rUQwdm9TTx8lhArNKbn8="
$wsCustomer.Range("I8").Value = "ASJ6u9M3Ltg5sLpg9eJmaD+P/czKtOBATxfgaZ8onQw="

$wsCustomer.Range("E2:E8").Value = $newDob

# ---------------------------------------------------------------
# Driver sheet - Password column J, rows 2-6
# ---------------------------------------------------------------
$wsDriver = $wb.Worksheets.Item("Driver")

$wsDriver.Range("J2").Value = "1VTZlApWK1DqezxcnW1fT4M+gtSeXhZyfP0MEkjLGd0="
$wsDriver.Range("J3").Value = "OJDmGlfoswwdTvJceIdcVxTJKRrB/YAvZe/6yz9ql7c="
$wsDriver.Range("J4").Value = "IqCkrdI9EwaAB494ALK2vfKkmzssvVZ2Oa+JeCYB/8k="
$wsDriver.Range("J5").Value = "dvfTgm52xqLEaCppWbQkKoSx0462FPNnE47rlxpktSo="
$wsDriver.Range("J6").Value = "gC+akQp0Bi2ll/kAbsIkBLaXiFmmThOcO7q+5TtoXI0="

$wsDriver.Range("E2:E6").Value = $newDob

# ---------------------------------------------------------------
# Owner sheet - Password column J, rows 2-4
# ---------------------------------------------------------------
$wsOwner = $wb.Worksheets.Item("Owner")

$wsOwner.Range("J2").Value = "DbJIZZBP74k7YyDX+v+d1aJxlc0vI4BNNqJClFbXTzA="
$wsOwner.Range("J3").Value = "XQK9M2z4I45FezKgmpYAdZb42DSBcaxzJQaR/1vsdkQ="
$wsOwner.Range("J4").Value = "6Vzzewqz54/eLvw78kfG79mV0fF3es1ljdOuta5WMzA="

$wsOwner.Range("E2:E4").Value = $newDob
